# Update cell values per the structural share index revision
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.02

# Row 3
$ws.Range("B3").Value = 0.05979281878947368
$ws.Range("C3").Value = 0.3735064921578948
$ws.Range("D3").Value = 0.19
$ws.Range("E3").Value = -0.3137136733684211

# Row 4
$ws.Range("B4").Value = 0.05262012604347826
$ws.Range("C4").Value = 0.3881498074347826
$ws.Range("D4").Value = 0.23
$ws.Range("E4").Value = -0.3355296813913043

# Row 6
$ws.Range("D6").Value = 0.02
$ws.Range("F6").Value = 0.02701135323073656

# Row 7
$ws.Range("B7").Value = 0.04913205644120356
$ws.Range("C7").Value = 0.3875504390053467
$ws.Range("D7").Value = 0.369
$ws.Range("E7").Value = -0.3384183825641432
$ws.Range("F7").Value = 0.07874922674061846

# Row 8
$ws.Range("B8").Value = 0.05796758252245629
$ws.Range("C8").Value = 0.3712761700058593
$ws.Range("D8").Value = 0.5945
$ws.Range("E8").Value = -0.313308587483403
$ws.Range("F8").Value = -0.06622690969025335

# Row 9
$ws.Range("B9").Value = 0.06742036616268397
$ws.Range("C9").Value = 0.3644014175974802
$ws.Range("D9").Value = 0.475
$ws.Range("E9").Value = -0.2969810514347962
$ws.Range("F9").Value = 0.2703176136005325

# Row 10
$ws.Range("B10").Value = 0.07938587643991264
$ws.Range("C10").Value = 0.3452797924278759
$ws.Range("D10").Value = 0.354
$ws.Range("E10").Value = -0.2658939159879632
$ws.Range("F10").Value = -0.1130255209680303

# Row 11
$ws.Range("B11").Value = 0.07766656881290657
$ws.Range("C11").Value = 0.3682778685654443
$ws.Range("D11").Value = 7.337000000000001
$ws.Range("E11").Value = -0.2906112997525377
$ws.Range("F11").Value = -0.1412662115142171

# Row 12
$ws.Range("B12").Value = 0.0761397651099326
$ws.Range("C12").Value = 0.3574056366816648
$ws.Range("D12").Value = 7.508
$ws.Range("E12").Value = -0.2812658715717322
$ws.Range("F12").Value = -0.1022720640026128

# Row 13
$ws.Range("B13").Value = 0.07597345355794274
$ws.Range("C13").Value = 0.3595361201589924
$ws.Range("D13").Value = 8.2025
$ws.Range("E13").Value = -0.2835626666010497
$ws.Range("F13").Value = -0.04518262956144392

# Row 14
$ws.Range("B14").Value = 0.07823492247708287
$ws.Range("C14").Value = 0.3549768740075087
$ws.Range("D14").Value = 8.3565
$ws.Range("E14").Value = -0.2767419515304259
$ws.Range("F14").Value = 0.04079835938387433

# Row 15
$ws.Range("B15").Value = 0.07778056904908576
$ws.Range("C15").Value = 0.3426512093598714
$ws.Range("D15").Value = 8.1675
$ws.Range("E15").Value = -0.2648706403107856
$ws.Range("F15").Value = -0.0885741864258921

# Row 16
$ws.Range("B16").Value = 0.07826656105167203
$ws.Range("C16").Value = 0.3374737219155959
$ws.Range("D16").Value = 7.9625
$ws.Range("E16").Value = -0.2592071608639239
$ws.Range("F16").Value = -0.07842654561871554

# Row 17
$ws.Range("B17").Value = 0.08048796327022274
$ws.Range("C17").Value = 0.3447974185939269
$ws.Range("D17").Value = 7.895
$ws.Range("E17").Value = -0.2643094553237042
$ws.Range("F17").Value = -0.06789755332789704

# Row 18
$ws.Range("B18").Value = 0.08110550422647941
$ws.Range("C18").Value = 0.3437312896671358
$ws.Range("D18").Value = 6.877000000000001
$ws.Range("E18").Value = -0.2626257854406564
$ws.Range("F18").Value = -0.05100840697156639

# Row 19
$ws.Range("B19").Value = 0.08336653195821253
$ws.Range("C19").Value = 0.3440032316670173
$ws.Range("D19").Value = 8.047000000000001
$ws.Range("E19").Value = -0.2606366997088048
$ws.Range("F19").Value = -0.01598493739061813

# Row 20
$ws.Range("B20").Value = 0.08383781642623232
$ws.Range("C20").Value = 0.3410041476302442
$ws.Range("D20").Value = 7.348
$ws.Range("E20").Value = -0.2571663312040119
$ws.Range("F20").Value = -0.007873353703308439

# Row 21
$ws.Range("B21").Value = 0.08395017232388889
$ws.Range("C21").Value = 0.3364210076848763
$ws.Range("D21").Value = 5.717000000000001
$ws.Range("E21").Value = -0.2524708353609874
$ws.Range("F21").Value = -0.0447907546410623

# Row 22
$ws.Range("B22").Value = 0.0822990903387351
$ws.Range("C22").Value = 0.3300620376686122
$ws.Range("D22").Value = 3.55
$ws.Range("E22").Value = -0.2477629473298771
$ws.Range("F22").Value = -0.05659321717340571

